# Add 6 extra "command" columns (X:AC) to the right of the last existing
# command column, pushing the old trailing column (the numeric total,
# previously in X) to AD. The new columns are filled with the same group
# value as the rest of the row (columns C:W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new columns at X:AC -- this shifts the old X:X column (and
# anything further right) to AD:AD, exactly like Excel's
# "Insert Copied/Entire Columns" shifting cells to the right.
$ws.Range("X1:AC1").EntireColumn.Insert()

# Fill the newly-inserted columns (X:AC) with each row's group value,
# matching the value already present in the neighbouring columns C:W.
$lastRow = 15
for ($r = 2; $r -le $lastRow; $r++) {
    $groupValue = $ws.Cells.Item($r, 23).Value2
    if ($groupValue -ne $null -and $groupValue -ne "") {
        $ws.Range("X" + $r + ":AC" + $r).Value = $groupValue
    }
}

Write-Output "Inserted 6 command columns (X:AC); shifted totals column to AD."
